$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Title for Stuart O'Neill (row 4) changes from "Mr" to "Dr."
$ws.Range("A4").Value = "Dr."

# Remove the hyperlink on the Stephen Jamison row before deleting it, so the
# relationship/hyperlink entry doesn't linger on a different cell after the
# row shift.
$ws.Range("E6").Hyperlinks.Delete()

# Drop the now-unused "Hyperlink" cell style definition.
try {
  $wb.Styles.Item("Hyperlink").Delete()
} catch {
}

# The whole Stephen Jamison row (row 6) is removed from the collaboration
# list for this draft circulation.
$ws.Rows("6:6").Delete()

# Reset view back to the top-left of the sheet with B3 selected (matches the
# draft's saved cursor position instead of the old F1-scrolled / H5-selected
# view).
$ws.Range("B3").Select()
